$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Weekly crime-stat table updates (rows 15-27) ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = 14.285714285714
$ws.Range("N15").Value = -55.555555555555
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 55.555555555555
$ws.Range("L16").Value = 21.739130434782
$ws.Range("M16").Value = -42.268041237113
$ws.Range("N16").Value = -79.562043795620
$ws.Range("C17").Value = 1
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 4
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 113
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = 1.801801801801
$ws.Range("L17").Value = 17.708333333333
$ws.Range("M17").Value = -9.6
$ws.Range("N17").Value = -59.786476868327
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 200
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 71
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 7.575757575757
$ws.Range("L18").Value = -12.345679012345
$ws.Range("M18").Value = -62.032085561497
$ws.Range("N18").Value = -94.151565074135
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 19.444444444444
$ws.Range("I19").Value = 320
$ws.Range("J19").Value = 294
$ws.Range("K19").Value = 8.843537414965
$ws.Range("L19").Value = 50.943396226415
$ws.Range("M19").Value = -13.043478260869
$ws.Range("N19").Value = -57.219251336898
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 61.538461538461
$ws.Range("I20").Value = 124
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = 74.647887323943
$ws.Range("L20").Value = 113.793103448276
$ws.Range("M20").Value = 18.095238095238
$ws.Range("N20").Value = -95.026073004412
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -38.095238095238
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 64
$ws.Range("H21").Value = 26.5625
$ws.Range("I21").Value = 693
$ws.Range("J21").Value = 587
$ws.Range("K21").Value = 18.057921635434
$ws.Range("L21").Value = 38.6
$ws.Range("M21").Value = -23
$ws.Range("N21").Value = -86.225402504472
$ws.Range("C23").Value = 2
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 19
$ws.Range("K23").Value = -45.714285714285
$ws.Range("L23").Value = 11.764705882352
$ws.Range("M23").Value = -9.523809523809
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 57.692307692307
$ws.Range("I24").Value = 770
$ws.Range("J24").Value = 446
$ws.Range("K24").Value = 72.645739910313
$ws.Range("L24").Value = 37.745974955277
$ws.Range("M24").Value = -49.308755760368
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 233.333333333333
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 76.470588235294
$ws.Range("I25").Value = 302
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = 36.651583710407
$ws.Range("L25").Value = 45.893719806763
$ws.Range("M25").Value = -37.474120082815
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = 100
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 18
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = 28.571428571428
$ws.Range("L26").Value = 28.571428571428
$ws.Range("C27").Value = 3
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 3
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = -3.333333333333
$ws.Range("L27").Value = 52.631578947368

Write-Output "edit applied"
